$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02300404818780313
$ws.Range("C2").Value = 0.337879676387987
$ws.Range("D2").Value = 0.2668330530619737
$ws.Range("E2").Value = 0.5165588573066711
$ws.Range("F2").Value = 0.5301872556910642
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.4362183256331637
$ws.Range("C3").Value = 0.5571722338028687
$ws.Range("D3").Value = 0.8814962611585107
$ws.Range("E3").Value = 0.9388803231288376
$ws.Range("F3").Value = 0.855493562520063
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.7297665574417316
$ws.Range("C4").Value = 0.7297665574417316
$ws.Range("D4").Value = 0.9114036334090323
$ws.Range("E4").Value = 0.9546746217476572
$ws.Range("F4").Value = 0.6344463573575141
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.6329849386531881
$ws.Range("C5").Value = 0.6568932063127731
$ws.Range("D5").Value = 0.6561129869449028
$ws.Range("E5").Value = 0.8100080165929858
$ws.Range("F5").Value = 0.5219890720522772
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.5673547796451167
$ws.Range("C6").Value = 0.5928185759272669
$ws.Range("D6").Value = 0.4654632102043935
$ws.Range("E6").Value = 0.6822486425082819
$ws.Range("F6").Value = 0.3922077131238274
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.5117157535595949
$ws.Range("C7").Value = 0.5117157535595949
$ws.Range("D7").Value = 0.3387636996991033
$ws.Range("E7").Value = 0.5820341052714207
$ws.Range("F7").Value = 0.287796619108349
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.4402895787565182
$ws.Range("C8").Value = 0.4402895787565182
$ws.Range("D8").Value = 0.2500964150522257
$ws.Range("E8").Value = 0.5000964057581555
$ws.Range("F8").Value = 0.2468365728875135
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.4177270777543842
$ws.Range("C9").Value = 0.4177270777543842
$ws.Range("D9").Value = 0.2295217317045278
$ws.Range("E9").Value = 0.4790842636786641
$ws.Range("F9").Value = 0.2488052405843257
$ws.Range("G9").Value = 9

$ws.Range("B10").Value = 0.2889056645135138
$ws.Range("C10").Value = 0.2889056645135138
$ws.Range("D10").Value = 0.1117927209483371
$ws.Range("E10").Value = 0.3343541848823446
$ws.Range("F10").Value = 0.188169597572051
$ws.Range("G10").Value = 5
